# Update faturamento_anual data for row 9 (Ano = 2025)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 4363314.1
$ws.Range("C9").Value = 698176.24
$ws.Range("D9").Value = 5061490.34
$ws.Range("E9").Value = 13.79388664406697
$ws.Range("F9").Value = 86.20611335593303
$ws.Range("G9").Value = -32.52474767122484
$ws.Range("H9").Value = -21.20457849959968
$ws.Range("I9").Value = 43665
$ws.Range("J9").Value = 1895
$ws.Range("K9").Value = 45560
$ws.Range("L9").Value = 31707
$ws.Range("M9").Value = 159.6332147475321
$ws.Range("N9").Value = 8.984847074262948
